# -----------------------------------------------------------------------
# Scheduled-runner update: refresh cached Universalis market-board pricing
# (currentAveragePrice* / LevePrice* / LeveProfit*) columns H:N for the
# leve rows whose item prices moved since the last data pull.
# -----------------------------------------------------------------------
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Leve Item ID 12604 (row 70)
$ws.Range("H70").Value = 5490.1
$ws.Range("I70").Value = 11400.25
$ws.Range("J70").Value = 1550
$ws.Range("K70").Value = 34200.75
$ws.Range("L70").Value = 4650
$ws.Range("M70").Value = -33930.75
$ws.Range("N70").Value = -5190

# Leve Item ID 12604 (row 73)
$ws.Range("H73").Value = 5490.1
$ws.Range("I73").Value = 11400.25
$ws.Range("J73").Value = 1550
$ws.Range("K73").Value = 34200.75
$ws.Range("L73").Value = 4650
$ws.Range("M73").Value = -33264.75
$ws.Range("N73").Value = -6522

# Leve Item ID 12603 (row 86)
$ws.Range("H86").Value = 2116.7666
$ws.Range("I86").Value = 1800.1875
$ws.Range("J86").Value = 2478.5715
$ws.Range("K86").Value = 1800.1875
$ws.Range("L86").Value = 2478.5715
$ws.Range("M86").Value = -677.1875
$ws.Range("N86").Value = -4724.5715

# Leve Item ID 12603 (row 89)
$ws.Range("H89").Value = 2116.7666
$ws.Range("I89").Value = 1800.1875
$ws.Range("J89").Value = 2478.5715
$ws.Range("K89").Value = 9000.9375
$ws.Range("L89").Value = 12392.8575
$ws.Range("M89").Value = -3384.9375
$ws.Range("N89").Value = -23624.8575

# Leve Item ID 19901 (row 92)
$ws.Range("H92").Value = 325.82758
$ws.Range("I92").Value = 315.3913
$ws.Range("K92").Value = 315.3913
$ws.Range("M92").Value = 932.6087

# Leve Item ID 34090 (row 123)
$ws.Range("H123").Value = 22500
$ws.Range("J123").Value = 22500
$ws.Range("L123").Value = 22500
$ws.Range("N123").Value = -32300

# Leve Item ID 44047 (row 135)
$ws.Range("H135").Value = 424.4
$ws.Range("I135").Value = 329.85715
$ws.Range("J135").Value = 1748
$ws.Range("K135").Value = 2968.71435
$ws.Range("L135").Value = 15732
$ws.Range("M135").Value = -433.7143499999997
$ws.Range("N135").Value = -20802

# Leve Item ID 44169 (row 138)
$ws.Range("H138").Value = 1368.1951
$ws.Range("I138").Value = 1152.4
$ws.Range("J138").Value = 10000
$ws.Range("K138").Value = 3457.2
$ws.Range("L138").Value = 30000
$ws.Range("M138").Value = 1682.8
$ws.Range("N138").Value = -40280

# Leve Item ID 44161 (row 141)
$ws.Range("H141").Value = 565.2368
$ws.Range("I141").Value = 555.1142599999999
$ws.Range("J141").Value = 683.3333
$ws.Range("K141").Value = 1665.34278
$ws.Range("L141").Value = 2049.9999
$ws.Range("M141").Value = 3514.65722
$ws.Range("N141").Value = -12409.9999


$ws = $wb.Worksheets.Item("ARM")
# Leve Item ID 43999 (row 61)
$ws.Range("H61").Value = 1569.6342
$ws.Range("I61").Value = 1097.6
$ws.Range("J61").Value = 2857
$ws.Range("K61").Value = 1097.6
$ws.Range("L61").Value = 2857
$ws.Range("M61").Value = -885.5999999999999
$ws.Range("N61").Value = -3281

# Leve Item ID 44000 (row 74)
$ws.Range("H74").Value = 3465.3333
$ws.Range("I74").Value = 4138.125
$ws.Range("J74").Value = 1312.4
$ws.Range("K74").Value = 4138.125
$ws.Range("L74").Value = 1312.4
$ws.Range("M74").Value = -3264.125
$ws.Range("N74").Value = -3060.4

# Leve Item ID 44000 (row 77)
$ws.Range("H77").Value = 3465.3333
$ws.Range("I77").Value = 4138.125
$ws.Range("J77").Value = 1312.4
$ws.Range("K77").Value = 20690.625
$ws.Range("L77").Value = 6562
$ws.Range("M77").Value = -16322.625
$ws.Range("N77").Value = -15298

# Leve Item ID 19941 (row 97)
$ws.Range("H97").Value = 501.65384
$ws.Range("I97").Value = 384.1579
$ws.Range("J97").Value = 820.5714
$ws.Range("K97").Value = 384.1579
$ws.Range("L97").Value = 820.5714
$ws.Range("M97").Value = 111.8421
$ws.Range("N97").Value = -1812.5714

# Leve Item ID 43999 (row 136)
$ws.Range("H136").Value = 1569.6342
$ws.Range("I136").Value = 1097.6
$ws.Range("J136").Value = 2857
$ws.Range("K136").Value = 3292.8
$ws.Range("L136").Value = 8571
$ws.Range("M136").Value = -742.7999999999997
$ws.Range("N136").Value = -13671


$ws = $wb.Worksheets.Item("BSM")
# Leve Item ID 43998 (row 134)
$ws.Range("H134").Value = 1500.3167
$ws.Range("I134").Value = 1289.3265
$ws.Range("J134").Value = 2440.182
$ws.Range("K134").Value = 3867.979499999999
$ws.Range("L134").Value = 7320.545999999999
$ws.Range("M134").Value = -1332.979499999999
$ws.Range("N134").Value = -12390.546


$ws = $wb.Worksheets.Item("CRP")
# Leve Item ID 44023 (row 31)
$ws.Range("H31").Value = 3409.262
$ws.Range("I31").Value = 2182.7646
$ws.Range("J31").Value = 4243.28
$ws.Range("K31").Value = 2182.7646
$ws.Range("L31").Value = 4243.28
$ws.Range("M31").Value = -1887.7646
$ws.Range("N31").Value = -4833.28

# Leve Item ID 44023 (row 34)
$ws.Range("H34").Value = 3409.262
$ws.Range("I34").Value = 2182.7646
$ws.Range("J34").Value = 4243.28
$ws.Range("K34").Value = 2182.7646
$ws.Range("L34").Value = 4243.28
$ws.Range("M34").Value = -1980.7646
$ws.Range("N34").Value = -4647.28


$ws = $wb.Worksheets.Item("CUL")
# Leve Item ID 43974 (row 5)
$ws.Range("H5").Value = 564.5484
$ws.Range("I5").Value = 496.26923
$ws.Range("J5").Value = 919.6
$ws.Range("K5").Value = 1488.80769
$ws.Range("L5").Value = 2758.8
$ws.Range("M5").Value = -1376.80769
$ws.Range("N5").Value = -2982.8

# Leve Item ID 12856 (row 82)
$ws.Range("H82").Value = 27500
$ws.Range("J82").Value = 27500
$ws.Range("L82").Value = 82500
$ws.Range("N82").Value = -83312

# Leve Item ID 12856 (row 85)
$ws.Range("H85").Value = 27500
$ws.Range("J85").Value = 27500
$ws.Range("L85").Value = 82500
$ws.Range("N85").Value = -85308

# Leve Item ID 43974 (row 135)
$ws.Range("H135").Value = 564.5484
$ws.Range("I135").Value = 496.26923
$ws.Range("J135").Value = 919.6
$ws.Range("K135").Value = 4466.42307
$ws.Range("L135").Value = 8276.4
$ws.Range("M135").Value = -1931.42307
$ws.Range("N135").Value = -13346.4


$ws = $wb.Worksheets.Item("GSM")
# Leve Item ID 19940 (row 97)
$ws.Range("H97").Value = 487.1111
$ws.Range("I97").Value = 504.875
$ws.Range("J97").Value = 345
$ws.Range("K97").Value = 504.875
$ws.Range("L97").Value = 345
$ws.Range("M97").Value = -8.875
$ws.Range("N97").Value = -1337

# Leve Item ID 44008 (row 132)
$ws.Range("H132").Value = 1927.4878
$ws.Range("I132").Value = 1774.5454
$ws.Range("J132").Value = 2558.375
$ws.Range("K132").Value = 5323.6362
$ws.Range("L132").Value = 7675.125
$ws.Range("M132").Value = -2793.6362
$ws.Range("N132").Value = -12735.125


$ws = $wb.Worksheets.Item("LTW")
# Leve Item ID 5289 (row 16)
$ws.Range("H16").Value = 745.63635
$ws.Range("I16").Value = 641.5714
$ws.Range("J16").Value = 927.75
$ws.Range("K16").Value = 641.5714
$ws.Range("L16").Value = 927.75
$ws.Range("M16").Value = -471.5714
$ws.Range("N16").Value = -1267.75

# Leve Item ID 12563 (row 68)
$ws.Range("H68").Value = 2999.8333
$ws.Range("I68").Value = 3000
$ws.Range("K68").Value = 3000
$ws.Range("M68").Value = -2251

# Leve Item ID 12563 (row 71)
$ws.Range("H71").Value = 2999.8333
$ws.Range("I71").Value = 3000
$ws.Range("K71").Value = 15000
$ws.Range("M71").Value = -11256

# Leve Item ID 10961 (row 88)
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("M88").ClearContents()

# Leve Item ID 10961 (row 91)
$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("M91").ClearContents()

# Leve Item ID 19993 (row 93)
$ws.Range("H93").Value = 8270
$ws.Range("I93").Value = 13402.75
$ws.Range("J93").Value = 1426.3334
$ws.Range("K93").Value = 13402.75
$ws.Range("L93").Value = 1426.3334
$ws.Range("M93").Value = -12154.75
$ws.Range("N93").Value = -3922.3334

# Leve Item ID 19995 (row 100)
$ws.Range("H100").Value = 63750724
$ws.Range("I100").Value = 1667465
$ws.Range("J100").Value = 250000500
$ws.Range("K100").Value = 1667465
$ws.Range("L100").Value = 250000500
$ws.Range("M100").Value = -1666924
$ws.Range("N100").Value = -250001582

# Leve Item ID 44058 (row 132)
$ws.Range("H132").Value = 1659.67
$ws.Range("I132").Value = 1659.4536
$ws.Range("J132").Value = 1666.6666
$ws.Range("K132").Value = 4978.3608
$ws.Range("L132").Value = 4999.9998
$ws.Range("M132").Value = -2448.3608
$ws.Range("N132").Value = -10059.9998

# Leve Item ID 44060 (row 136)
$ws.Range("H136").Value = 1895.1267
$ws.Range("I136").Value = 1487.5741
$ws.Range("J136").Value = 3189.7058
$ws.Range("K136").Value = 4462.7223
$ws.Range("L136").Value = 9569.117400000001
$ws.Range("M136").Value = -1912.7223
$ws.Range("N136").Value = -14669.1174


$ws = $wb.Worksheets.Item("WVR")
# Leve Item ID 19977 (row 96)
$ws.Range("H96").Value = 8103.067
$ws.Range("I96").Value = 1024
$ws.Range("J96").Value = 10677.272
$ws.Range("K96").Value = 1024
$ws.Range("L96").Value = 10677.272
$ws.Range("M96").Value = 349
$ws.Range("N96").Value = -13423.272

# Leve Item ID 19981 (row 100)
$ws.Range("H100").Value = 740.4
$ws.Range("I100").Value = 617.3333
$ws.Range("J100").Value = 925
$ws.Range("K100").Value = 1234.6666
$ws.Range("L100").Value = 1850
$ws.Range("M100").Value = -693.6666
$ws.Range("N100").Value = -2932

# Leve Item ID 44031 (row 136)
$ws.Range("H136").Value = 2446.4546
$ws.Range("I136").Value = 2333.8333
$ws.Range("J136").Value = 2953.25
$ws.Range("K136").Value = 7001.499899999999
$ws.Range("L136").Value = 8859.75
$ws.Range("M136").Value = -4451.499899999999
$ws.Range("N136").Value = -13959.75
